$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '#leand'
$ws.Range("C2").Value = 'Leand'
$ws.Range("D2").ClearContents()
$ws.Range("B3").Value = '#diew'
$ws.Range("C3").Value = 'Diew'
$ws.Range("D3").ClearContents()
$ws.Range("B4").Value = '#machtelt-met-haer-kint-op-d''arm'
$ws.Range("C4").Value = 'Machtelt met haer kint op d''arm'
$ws.Range("D4").ClearContents()
$ws.Range("B5").Value = '#hans'
$ws.Range("C5").Value = 'Hans'
$ws.Range("D5").ClearContents()
$ws.Range("B6").Value = '#leander-singht-twee-veersjens'
$ws.Range("C6").Value = 'Leander singht twee veersjens'
$ws.Range("D6").ClearContents()
$ws.Range("B7").Value = '#leander-weêr-uyt'
$ws.Range("C7").Value = 'Leander weêr uyt'
$ws.Range("D7").ClearContents()
$ws.Range("B8").Value = '#bea'
$ws.Range("C8").Value = 'Bea'
$ws.Range("D8").ClearContents()
$ws.Range("B9").Value = '#iurien'
$ws.Range("C9").Value = 'Iurien'
$ws.Range("D9").ClearContents()
$ws.Range("B10").Value = '#beatrix'
$ws.Range("C10").Value = 'Beatrix'
$ws.Range("D10").ClearContents()
$ws.Range("B11").Value = '#leander-uyt'
$ws.Range("C11").Value = 'Leander uyt'
$ws.Range("D11").ClearContents()
$ws.Range("B12").Value = '#diewer'
$ws.Range("C12").Value = 'Diewer'
$ws.Range("D12").ClearContents()
$ws.Range("B13").Value = '#pleuntje-ians'
$ws.Range("C13").Value = 'Pleuntje Ians'
$ws.Range("D13").ClearContents()
$ws.Range("B14").Value = '#grietje'
$ws.Range("C14").Value = 'Grietje'
$ws.Range("D14").ClearContents()
$ws.Range("B15").Value = '#pleuntje'
$ws.Range("C15").Value = 'Pleuntje'
$ws.Range("D15").ClearContents()
$ws.Range("B16").Value = '#anne-pieters'
$ws.Range("C16").Value = 'Anne Pieters'
$ws.Range("D16").ClearContents()
$ws.Range("B17").Value = '#ioost'
$ws.Range("C17").Value = 'Ioost'
$ws.Range("D17").ClearContents()
$ws.Range("B18").Value = '#anne-pieters.-pleuntje-ians'
$ws.Range("C18").Value = 'Anne Pieters. Pleuntje Ians'
$ws.Range("D18").ClearContents()
$ws.Range("B19").Value = '#beat'
$ws.Range("C19").Value = 'Beat'
$ws.Range("D19").ClearContents()
$ws.Range("B20").Value = '#macht'
$ws.Range("C20").Value = 'Macht'
$ws.Range("D20").ClearContents()
$ws.Range("B21").Value = '#anne'
$ws.Range("C21").Value = 'Anne'
$ws.Range("D21").ClearContents()
$ws.Range("B22").Value = '#iochem'
$ws.Range("C22").Value = 'Iochem'
$ws.Range("D22").ClearContents()
$ws.Range("B23").Value = '#leander.-anne-pieters.-hans.-iurien.-iochgum'
$ws.Range("C23").Value = 'Leander. Anne Pieters. Hans. Iurien. Iochgum'
$ws.Range("D23").ClearContents()
$ws.Range("B24").Value = '#leander'
$ws.Range("C24").Value = 'Leander'
$ws.Range("D24").ClearContents()
$ws.Range("B25").Value = '#grietje-gaet-aen-''t-loogen'
$ws.Range("C25").Value = 'Grietje gaet aen ''t loogen'
$ws.Range("D25").ClearContents()
$ws.Range("B26").Value = '#grietje-met-een-hallefvat'
$ws.Range("C26").Value = 'Grietje met een hallefvat'
$ws.Range("D26").ClearContents()
$ws.Range("B27").Value = '#grietje.-leander'
$ws.Range("C27").Value = 'Grietje. Leander'
$ws.Range("D27").ClearContents()
